$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.569.36"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.115.80"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.15%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.91%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "337.05"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.93%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.58%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +4.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.37"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.93%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.23%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.82%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.56"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.116.52"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.871"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.119"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +5.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001173"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +4.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "97.07"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.010"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06686"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.35%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.42"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.314"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.646.22"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.361"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.73%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.369.21"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.35"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "164.08"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.547"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.38%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "134.21"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.215"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1072"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.651"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.377"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.46%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.951"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.37%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.60"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.896"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +7.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02626"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06846"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2329"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.35%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6888"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.260"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.70%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.90"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +7.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6460"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.325"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.94%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000367"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +23.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.686"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.67%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.256"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "83.51"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3342"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +11.62%  "
